$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false
$d = $word.ActiveDocument
$changeCount = 0

# P4 run 4
$para = $d.Paragraphs.Item(4)
$rng = $para.Range
$ok = $rng.Find.Execute("Ativação: 01/01/2025", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 15/07/2025", 2)
if (-not $ok) { Write-Output "FAILED P4 run 4" } else { $changeCount++ }

# P6 run 1
$para = $d.Paragraphs.Item(6)
$rng = $para.Range
$ok = $rng.Find.Execute("Desenvolver habilidades avançadas de comunicação oral e escrita, concentrando-se na divulgação científica em biotecnologia. Capacitar os alunos para comunicar conceitos complexos de biotecnologia de forma acessível ao público leigo. Promover a interação entre ciência e sociedade, estimulando a reflexão sobre o impacto social, ético e econômico da biotecnologia. Fomentar uma participação cidadã informada, visando combater a disseminação de desinformação. Estimular a criatividade e a inovação na comunicação científica, incentivando abordagens inovadoras e explorando diversas mídias sociais para alcançar públicos variados.", $true, $false, $false, $false, $false, $true, 1, $false, "1. O método científico. 2. Introdução à Biotecnologia. 3. Importância da Divulgação Científica. 4. Ética na Comunicação Científica. 5. Mídias sociais na Divulgação Científica. 6. Escrita Científica. 7. Visita supervisionada prevista.", 2)
if (-not $ok) { Write-Output "FAILED P6 run 1" } else { $changeCount++ }

# P7 run 1
$para = $d.Paragraphs.Item(7)
$rng = $para.Range
$ok = $rng.Find.Execute("To develop advanced oral and written communication skills, focusing on scientific dissemination in biotechnology. Enable students on how to communicate complex biotechnology concepts in an accessible manner to the general public. Foster interaction between science and society, encouraging reflection on the social, ethical, and economic impact of biotechnology. Promote informed civic participation, aiming to combat the spread of misinformation. Stimulate creativity and innovation in scientific communication, encouraging innovative approaches and exploring various social medias to reach diverse audiences.", $true, $false, $false, $false, $false, $true, 1, $false, "1.Introduction to Biotechnology. 2. Importance of Scientific Dissemination. 3. Ethics in Scientific Communication. 4. Social media in Scientific Outreach. 5. Scientific Writing. 6. Supervised visits.", 2)
if (-not $ok) { Write-Output "FAILED P7 run 1" } else { $changeCount++ }

# P9 run 2
$para = $d.Paragraphs.Item(9)
$rng = $para.Range
$ok = $rng.Find.Execute("5111420 - Talita Martins Lacerda", $true, $false, $false, $false, $false, $true, 1, $false, "1. Introdução sobre os diversos tipos de conhecimento e dos métodos científicos e sobre as principais técnicas utilizadas para coleta de dados que permitem o desenvolvimento das pesquisas científicas;^l2. Introdução à Biotecnologia: Aplicações práticas da biotecnologia nos diversos campos científicos. Porque divulgá-los?^l3. Importância da Divulgação Científica: Contextualização histórica. Impacto na sociedade. Estratégias de comunicação. Técnicas de simplificação - Adaptação de linguagem para diferentes públicos.^l4. Ética na Comunicação Científica: Responsabilidade social. Transparência e honestidade. Abordagem da era das `"Fake News`" e seus desdobramentos.^l5. Mídias sociais na Divulgação Científica: Vídeos, podcasts e posts. Exploração prática de plataformas como Facebook, Instagram, LinkedIn, ResearchGate, entre outras.^l6. Escrita Científica: Artigos científicos. Métricas científicas (fator de impacto, índice H, etc.). Tipos de artigos e estratégias de publicação. Relatórios e outros documentos técnicos.^l7. Atividades Práticas: Criação de materiais de comunicação científica inovadores e eficazes em biotecnologia (elaboração de apresentações, infográficos e vídeos de divulgação científica).^l8. Visitas supervisionadas a laboratórios e indústrias (viagem didática complementar), a depender da viabilidade no momento do oferecimento da disciplina.", 2)
if (-not $ok) { Write-Output "FAILED P9 run 2" } else { $changeCount++ }

# P9 run 1
$para = $d.Paragraphs.Item(9)
$rng = $para.Range
$ok = $rng.Find.Execute("5082401 - André Moreni Lopes", $true, $false, $false, $false, $false, $true, 1, $false, "Desenvolver habilidades avançadas de comunicação oral e escrita, concentrando-se na divulgação científica em biotecnologia. Capacitar os alunos para comunicar conceitos complexos de biotecnologia de forma acessível ao público leigo. Promover a interação entre ciência e sociedade, estimulando a reflexão sobre o impacto social, ético e econômico da biotecnologia. Fomentar uma participação cidadã informada, visando combater a disseminação de desinformação. Estimular a criatividade e a inovação na comunicação científica, incentivando abordagens inovadoras e explorando diversas mídias sociais para alcançar públicos variados.", 2)
if (-not $ok) { Write-Output "FAILED P9 run 1" } else { $changeCount++ }

# P11 run 1
$para = $d.Paragraphs.Item(11)
$rng = $para.Range
$ok = $rng.Find.Execute("1. O método científico. 2. Introdução à Biotecnologia. 3. Importância da Divulgação Científica. 4. Ética na Comunicação Científica. 5. Mídias sociais na Divulgação Científica. 6. Escrita Científica.", $true, $false, $false, $false, $false, $true, 1, $false, "Notas `"N`" distribuídas da seguinte forma:^lParticipação em Atividades Práticas (N = 30%), Projetos de Divulgação Científica (N = 40%), Avaliação Teórica Final (N = 20%), Contribuição para Eventos de Extensão Universitária (N = 10%).", 2)
if (-not $ok) { Write-Output "FAILED P11 run 1" } else { $changeCount++ }

# P12 run 1
$para = $d.Paragraphs.Item(12)
$rng = $para.Range
$ok = $rng.Find.Execute("1.Introduction to Biotechnology. 2. Importance of Scientific Dissemination. 3. Ethics in Scientific Communication. 4. Social media in Scientific Outreach. 5. Scientific Writing.", $true, $false, $false, $false, $false, $true, 1, $false, "To develop advanced oral and written communication skills, focusing on scientific dissemination in biotechnology. Enable students on how to communicate complex biotechnology concepts in an accessible manner to the general public. Foster interaction between science and society, encouraging reflection on the social, ethical, and economic impact of biotechnology. Promote informed civic participation, aiming to combat the spread of misinformation. Stimulate creativity and innovation in scientific communication, encouraging innovative approaches and exploring various social medias to reach diverse audiences.", 2)
if (-not $ok) { Write-Output "FAILED P12 run 1" } else { $changeCount++ }

# P14 run 1
$para = $d.Paragraphs.Item(14)
$rng = $para.Range
$ok = $rng.Find.Execute("1. Introdução sobre os diversos tipos de conhecimento e dos métodos científicos e sobre as principais técnicas utilizadas para coleta de dados que permitem o desenvolvimento das pesquisas científicas;^l2. Introdução à Biotecnologia: Aplicações práticas da biotecnologia nos diversos campos científicos. Porque divulgá-los?^l3. Importância da Divulgação Científica: Contextualização histórica. Impacto na sociedade. Estratégias de comunicação. Técnicas de simplificação - Adaptação de linguagem para diferentes públicos.^l4. Ética na Comunicação Científica: Responsabilidade social. Transparência e honestidade. Abordagem da era das `"Fake News`" e seus desdobramentos.^l5. Mídias sociais na Divulgação Científica: Vídeos, podcasts e posts. Exploração prática de plataformas como Facebook, Instagram, LinkedIn, ResearchGate, entre outras.^l6. Escrita Científica: Artigos científicos. Métricas científicas (fator de impacto, índice H, etc.). Tipos de artigos e estratégias de publicação. Relatórios e outros documentos técnicos.^l7. Atividades Práticas: Criação de materiais de comunicação científica inovadores e eficazes em biotecnologia (elaboração de apresentações, infográficos e vídeos de divulgação científica).", $true, $false, $false, $false, $false, $true, 1, $false, "Os critérios de avaliação serão definidos e informados pelo(s) docente(s) responsável(is) em momento oportuno, e a Média Final (MF) será calculada considerando-se o peso de cada N.", 2)
if (-not $ok) { Write-Output "FAILED P14 run 1" } else { $changeCount++ }

# P15 run 1
$para = $d.Paragraphs.Item(15)
$rng = $para.Range
$ok = $rng.Find.Execute("1. Introduction to the various types of knowledge and to the scientific methods, as well as the main techniques used for data collection that enable the development of scientific research;^l2. Introduction to Biotechnology: Practical applications of biotechnology in various scientific fields. Why to talk about them with society?^l3. Importance of Scientific Dissemination: Historical context. Impact on society. Communication strategies. Simplification techniques - Adapting language for different audiences;^l4. Ethics in Scientific Communication: Social responsibility. Transparency and honesty. Addressing the era of `"Fake News`" and its consequences;^l5. Social media in Scientific Dissemination: Videos, podcasts, and posts. Practical exploration of platforms such as Facebook, Instagram, LinkedIn, ResearchGate, among others;^l6. Scientific Writing: Scientific articles. Scientific metrics (impact factor, H-index, etc.). Types of articles and publication strategies. Reports and other technical documents;^l7. Practical Activities: Creation of innovative and effective scientific communication materials in biotechnology (development of presentations, infographics, and scientific dissemination videos).", $true, $false, $false, $false, $false, $true, 1, $false, "1. Introduction to the various types of knowledge and to the scientific methods, as well as the main techniques used for data collection that enable the development of scientific research;^l2. Introduction to Biotechnology: Practical applications of biotechnology in various scientific fields. Why to talk about them with society?^l3. Importance of Scientific Dissemination: Historical context. Impact on society. Communication strategies. Simplification techniques - Adapting language for different audiences;^l4. Ethics in Scientific Communication: Social responsibility. Transparency and honesty. Addressing the era of `"Fake News`" and its consequences;^l5. Social media in Scientific Dissemination: Videos, podcasts, and posts. Practical exploration of platforms such as Facebook, Instagram, LinkedIn, ResearchGate, among others;^l6. Scientific Writing: Scientific articles. Scientific metrics (impact factor, H-index, etc.). Types of articles and publication strategies. Reports and other technical documents;^l7. Practical Activities: Creation of innovative and effective scientific communication materials in biotechnology (development of presentations, infographics, and scientific dissemination videos).^l8. Supervised visits to laboratories and industries, depending on feasibility at the time the discipline is offered", 2)
if (-not $ok) { Write-Output "FAILED P15 run 1" } else { $changeCount++ }

# P17 run 6
$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$ok = $rng.Find.Execute("Uma Prova de Recuperação (PR) será aplicada para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Serão considerados aprovados os alunos que tenham obtido Nota Final (NF) igual ou maior do que 5,0.", $true, $false, $false, $false, $false, $true, 1, $false, "5082401 - André Moreni Lopes", 2)
if (-not $ok) { Write-Output "FAILED P17 run 6" } else { $changeCount++ }

# P17 run 4
$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$ok = $rng.Find.Execute("Os critérios de avaliação serão definidos e informados pelo(s) docente(s) responsável(is) em momento oportuno, e a Média Final (MF) será calculada considerando-se o peso de cada N.", $true, $false, $false, $false, $false, $true, 1, $false, "-SEVERINO, A. J. Metodologia do trabalho científico [livro eletrônico]. 1ª Ed. -- São Paulo: Cortez, 2013.^l-VOGT, C., GOMES, M., MUNIZ, R. (2018). ComCiência e divulgação científica.^l-PORTO, C., BROTAS, A., BORTOLIERO, S. (2011). Diálogos entre ciência e divulgação científica: leituras contemporâneas. EDUFBA.^l-BUCCHI, M., TRENCH, B. (Eds.). (2021). Routledge handbook of public communication of science and technology. Routledge.^l-CHALMERS, A.F. O que é ciência afinal? Trad. Raul Filker. São Paulo: Editora Brasiliense, 1993.^l-SEVERINO, A. J. 2007. Metodologia do trabalho científico. 23a Ed. revista e atualizada. Ed. Cortez, São Paulo.^l-GALLIANO, A. G. O método científico: teoria e prática. São Paulo: Harbra, 1986.^l-SEVERINO, Antônio Joaquim. Metodologia do Trabalho Científico. São Paulo: CORTEZ, 2008.^l-Redação Científica: http://www.gilsonvolpato.com.br/", 2)
if (-not $ok) { Write-Output "FAILED P17 run 4" } else { $changeCount++ }

# P17 run 2
$para = $d.Paragraphs.Item(17)
$rng = $para.Range
$ok = $rng.Find.Execute("Notas `"N`" distribuídas da seguinte forma:^lParticipação em Atividades Práticas (N = 30%), Projetos de Divulgação Científica (N = 40%), Avaliação Teórica Final (N = 20%), Contribuição para Eventos de Extensão Universitária (N = 10%).", $true, $false, $false, $false, $false, $true, 1, $false, "Uma Prova de Recuperação (PR) será aplicada para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Serão considerados aprovados os alunos que tenham obtido Nota Final (NF) igual ou maior do que 5,0.", 2)
if (-not $ok) { Write-Output "FAILED P17 run 2" } else { $changeCount++ }

# P19 run 1
$para = $d.Paragraphs.Item(19)
$rng = $para.Range
$ok = $rng.Find.Execute("-SEVERINO, A. J. Metodologia do trabalho científico [livro eletrônico]. 1ª Ed. -- São Paulo: Cortez, 2013.^l-VOGT, C., GOMES, M., MUNIZ, R. (2018). ComCiência e divulgação científica.^l-PORTO, C., BROTAS, A., BORTOLIERO, S. (2011). Diálogos entre ciência e divulgação científica: leituras contemporâneas. EDUFBA.^l-BUCCHI, M., TRENCH, B. (Eds.). (2021). Routledge handbook of public communication of science and technology. Routledge.^l-CHALMERS, A.F. O que é ciência afinal? Trad. Raul Filker. São Paulo: Editora Brasiliense, 1993.^l-SEVERINO, A. J. 2007. Metodologia do trabalho científico. 23a Ed. revista e atualizada. Ed. Cortez, São Paulo.^l-GALLIANO, A. G. O método científico: teoria e prática. São Paulo: Harbra, 1986.^l-SEVERINO, Antônio Joaquim. Metodologia do Trabalho Científico. São Paulo: CORTEZ, 2008.^l-Redação Científica: http://www.gilsonvolpato.com.br/", $true, $false, $false, $false, $false, $true, 1, $false, "5111420 - Talita Martins Lacerda", 2)
if (-not $ok) { Write-Output "FAILED P19 run 1" } else { $changeCount++ }

Write-Output "Total replacements applied: $changeCount"